$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investment_Cost")

# Update row 14 values (Distillation_tower_FT investment costs)
$ws.Range("B14").Value = 2233000
$ws.Range("C14").Value = 1967000
$ws.Range("D14").Value = 1701000
$ws.Range("E14").Value = 1170000
$ws.Range("F14").Value = 957000

# E14 adopts the same formatting as F14 (General number format, top-aligned wrap text)
$ws.Range("F14").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the note in I14 to the new note text
$ws.Range("I14").Value = "No value given for 2025 (linear approximated)"

# Update the view selection (the sheet is no longer scrolled to D6; C20 is now selected)
$ws.Activate()
$ws.Range("C20").Select()
